$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = 44362  # D2
$ws.Cells.Item(2, 10).Value = 25  # J2
$ws.Cells.Item(2, 11).Value = 8000  # K2
$ws.Cells.Item(2, 12).Value = 8000  # L2
$ws.Cells.Item(2, 13).Value = 8000  # M2
$ws.Cells.Item(2, 15).Value = 'Región Metropolitana'  # O2
$ws.Cells.Item(2, 16).Value = 500  # P2
$ws.Cells.Item(3, 4).Value = 44354  # D3
$ws.Cells.Item(3, 10).Value = 100  # J3
$ws.Cells.Item(3, 11).Value = 8000  # K3
$ws.Cells.Item(3, 12).Value = 9000  # L3
$ws.Cells.Item(3, 13).Value = 8500  # M3
$ws.Cells.Item(3, 16).Value = 531  # P3
$ws.Cells.Item(4, 4).Value = 44354  # D4
$ws.Cells.Item(4, 10).Value = 80  # J4
$ws.Cells.Item(4, 11).Value = 9000  # K4
$ws.Cells.Item(4, 12).Value = 9000  # L4
$ws.Cells.Item(4, 13).Value = 9000  # M4
$ws.Cells.Item(4, 16).Value = 562  # P4
$ws.Cells.Item(5, 4).Value = 44355  # D5
$ws.Cells.Item(5, 10).Value = 30  # J5
$ws.Cells.Item(5, 11).Value = 8000  # K5
$ws.Cells.Item(5, 12).Value = 8000  # L5
$ws.Cells.Item(5, 13).Value = 8000  # M5
$ws.Cells.Item(5, 14).Value = '$/caja 16 unidades'  # N5
$ws.Cells.Item(5, 16).Value = 500  # P5
$ws.Cells.Item(5, 17).Value = 16  # Q5
$ws.Cells.Item(6, 4).Value = 44389  # D6
$ws.Cells.Item(6, 10).Value = 55  # J6
$ws.Cells.Item(7, 4).Value = 44305  # D7
$ws.Cells.Item(7, 10).Value = 35  # J7
$ws.Cells.Item(8, 4).Value = 44312  # D8
$ws.Cells.Item(8, 10).Value = 40  # J8
$ws.Cells.Item(8, 11).Value = 7000  # K8
$ws.Cells.Item(8, 12).Value = 7000  # L8
$ws.Cells.Item(8, 13).Value = 7000  # M8
$ws.Cells.Item(8, 15).Value = 'Región del Maule'  # O8
$ws.Cells.Item(8, 16).Value = 438  # P8
$ws.Cells.Item(9, 4).Value = 44403  # D9
$ws.Cells.Item(9, 10).Value = 35  # J9
$ws.Cells.Item(9, 11).Value = 5000  # K9
$ws.Cells.Item(9, 12).Value = 5000  # L9
$ws.Cells.Item(9, 13).Value = 5000  # M9
$ws.Cells.Item(9, 15).Value = 'Región Metropolitana'  # O9
$ws.Cells.Item(9, 16).Value = 312  # P9
$ws.Cells.Item(10, 4).Value = 44676  # D10
$ws.Cells.Item(10, 9).Value = 'Primera'  # I10
$ws.Cells.Item(10, 10).Value = 40  # J10
$ws.Cells.Item(10, 11).Value = 12000  # K10
$ws.Cells.Item(10, 12).Value = 12000  # L10
$ws.Cells.Item(10, 13).Value = 12000  # M10
$ws.Cells.Item(10, 14).Value = '$/caja 18 unidades'  # N10
$ws.Cells.Item(10, 15).Value = 'Región Metropolitana'  # O10
$ws.Cells.Item(10, 16).Value = 667  # P10
$ws.Cells.Item(10, 17).Value = 18  # Q10
$ws.Cells.Item(11, 4).Value = 44372  # D11
$ws.Cells.Item(11, 10).Value = 50  # J11
$ws.Cells.Item(11, 11).Value = 6000  # K11
$ws.Cells.Item(11, 13).Value = 6400  # M11
$ws.Cells.Item(11, 16).Value = 400  # P11
$ws.Cells.Item(12, 4).Value = 44396  # D12
$ws.Cells.Item(12, 10).Value = 80  # J12
$ws.Cells.Item(14, 4).Value = 44420  # D14
$ws.Cells.Item(14, 10).Value = 45  # J14
$ws.Cells.Item(15, 4).Value = 44313  # D15
$ws.Cells.Item(15, 10).Value = 20  # J15
$ws.Cells.Item(15, 11).Value = 7000  # K15
$ws.Cells.Item(15, 12).Value = 7000  # L15
$ws.Cells.Item(15, 13).Value = 7000  # M15
$ws.Cells.Item(15, 15).Value = 'Región del Maule'  # O15
$ws.Cells.Item(15, 16).Value = 438  # P15
$ws.Cells.Item(16, 4).Value = 44392  # D16
$ws.Cells.Item(16, 10).Value = 95  # J16
$ws.Cells.Item(16, 11).Value = 7000  # K16
$ws.Cells.Item(16, 12).Value = 7000  # L16
$ws.Cells.Item(16, 13).Value = 7000  # M16
$ws.Cells.Item(16, 15).Value = 'Región del Maule'  # O16
$ws.Cells.Item(16, 16).Value = 438  # P16
$ws.Cells.Item(17, 4).Value = 44308  # D17
$ws.Cells.Item(17, 10).Value = 75  # J17
$ws.Cells.Item(17, 11).Value = 5000  # K17
$ws.Cells.Item(17, 12).Value = 5000  # L17
$ws.Cells.Item(17, 13).Value = 5000  # M17
$ws.Cells.Item(17, 16).Value = 312  # P17
$ws.Cells.Item(18, 4).Value = 44386  # D18
$ws.Cells.Item(18, 10).Value = 40  # J18
$ws.Cells.Item(18, 11).Value = 7000  # K18
$ws.Cells.Item(18, 12).Value = 7000  # L18
$ws.Cells.Item(18, 13).Value = 7000  # M18
$ws.Cells.Item(18, 15).Value = 'Región del Maule'  # O18
$ws.Cells.Item(18, 16).Value = 438  # P18
$ws.Cells.Item(19, 4).Value = 44371  # D19
$ws.Cells.Item(19, 10).Value = 200  # J19
$ws.Cells.Item(19, 15).Value = 'Región Metropolitana'  # O19
$ws.Cells.Item(20, 4).Value = 44385  # D20
$ws.Cells.Item(20, 10).Value = 100  # J20
$ws.Cells.Item(20, 11).Value = 7000  # K20
$ws.Cells.Item(20, 12).Value = 7000  # L20
$ws.Cells.Item(20, 13).Value = 7000  # M20
$ws.Cells.Item(20, 15).Value = 'Región del Maule'  # O20
$ws.Cells.Item(20, 16).Value = 438  # P20
$ws.Cells.Item(21, 4).Value = 44397  # D21
$ws.Cells.Item(21, 10).Value = 40  # J21
$ws.Cells.Item(21, 11).Value = 8000  # K21
$ws.Cells.Item(21, 12).Value = 8000  # L21
$ws.Cells.Item(21, 13).Value = 8000  # M21
$ws.Cells.Item(21, 15).Value = 'Región Metropolitana'  # O21
$ws.Cells.Item(21, 16).Value = 500  # P21
$ws.Cells.Item(22, 4).Value = 44348  # D22
$ws.Cells.Item(22, 10).Value = 35  # J22
$ws.Cells.Item(23, 4).Value = 44315  # D23
$ws.Cells.Item(23, 10).Value = 40  # J23
$ws.Cells.Item(23, 15).Value = 'Región del Maule'  # O23
$ws.Cells.Item(24, 4).Value = 44398  # D24
$ws.Cells.Item(24, 10).Value = 80  # J24
$ws.Cells.Item(24, 15).Value = 'Región Metropolitana'  # O24
$ws.Cells.Item(25, 4).Value = 44314  # D25
$ws.Cells.Item(25, 9).Value = 'Segunda'  # I25
$ws.Cells.Item(25, 10).Value = 20  # J25
$ws.Cells.Item(25, 11).Value = 5000  # K25
$ws.Cells.Item(25, 12).Value = 5000  # L25
$ws.Cells.Item(25, 13).Value = 5000  # M25
$ws.Cells.Item(25, 16).Value = 312  # P25
$ws.Cells.Item(26, 4).Value = 44467  # D26
$ws.Cells.Item(26, 10).Value = 40  # J26
$ws.Cells.Item(26, 15).Value = 'Región del Maule'  # O26
$ws.Cells.Item(27, 4).Value = 44369  # D27
$ws.Cells.Item(27, 10).Value = 60  # J27
$ws.Cells.Item(27, 15).Value = 'Región Metropolitana'  # O27
